$wb = $excel.ActiveWorkbook

# New client row to insert (alphabetically between "FRANK FERRETERIA..." and
# "ILLER LOPEZ ROBERTO FERNANDO") for advisor "ILLER LOPEZ ROBERTO FERNANDO".
$asesor = "ILLER LOPEZ ROBERTO FERNANDO"
$nuevoCliente = "HUERTA MUÑOZ NANCY ELIZABETH"
$insertRow = 11
$newTotalRow = 19      # the moved-down summary/total row after the insert

# ---- Sheet "VENTAS POR GRUPO" (columns A:R) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$lastCol1 = 18   # column R

$ws1.Rows.Item($insertRow).Insert()
$ws1.Cells.Item($insertRow, 1).Value = $asesor
$ws1.Cells.Item($insertRow, 2).Value = $nuevoCliente
for ($c = 3; $c -le $lastCol1; $c++) {
    $ws1.Cells.Item($insertRow, $c).Value = 0
}

# Update the "X de 16" -> "X de 17" summary labels on the (now shifted) total row.
for ($c = 3; $c -le $lastCol1; $c++) {
    $cell = $ws1.Cells.Item($newTotalRow, $c)
    $txt = $cell.Text
    if ($txt -match "^(\d+) de \d+$") {
        $cell.Value = "$($Matches[1]) de 17"
    }
}

# ---- Sheet "VENTA MENSUAL" (columns A:G) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$lastCol2 = 7   # column G

$ws2.Rows.Item($insertRow).Insert()
$ws2.Cells.Item($insertRow, 1).Value = $asesor
$ws2.Cells.Item($insertRow, 2).Value = $nuevoCliente
for ($c = 3; $c -le $lastCol2; $c++) {
    $ws2.Cells.Item($insertRow, $c).Value = 0
}
